$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at row 34, pushing the existing rows 34-37 down to 36-39
# (their contents remain unchanged by the shift).
$ws.Rows.Item(34).Resize(2).Insert()

# New row 34: latest weekly record (same shape as the old row 34, new date)
$ws.Range("A34").Value = 7
$ws.Range("B34").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C34").Value = "Ñuble"
$ws.Range("D34").Value = 45244
$ws.Range("E34").Value = 16
$ws.Range("F34").Value = 100112039
$ws.Range("G34").Value = "Ciboulette"
$ws.Range("H34").Value = "Sin especificar"
$ws.Range("I34").Value = "Primera"
$ws.Range("J34").Value = 100
$ws.Range("K34").Value = 2500
$ws.Range("L34").Value = 2500
$ws.Range("M34").Value = 2500
$ws.Range("N34").Value = "$/docena de atados"
$ws.Range("O34").Value = "Región Metropolitana"
$ws.Range("P34").Value = 833
$ws.Range("Q34").Value = 3
$ws.Range("R34").Value = "Hortaliza"

# New row 35: additional weekly record (Segunda quality)
$ws.Range("A35").Value = 7
$ws.Range("B35").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C35").Value = "Ñuble"
$ws.Range("D35").Value = 45244
$ws.Range("E35").Value = 16
$ws.Range("F35").Value = 100112039
$ws.Range("G35").Value = "Ciboulette"
$ws.Range("H35").Value = "Sin especificar"
$ws.Range("I35").Value = "Segunda"
$ws.Range("J35").Value = 100
$ws.Range("K35").Value = 2000
$ws.Range("L35").Value = 2000
$ws.Range("M35").Value = 2000
$ws.Range("N35").Value = "$/docena de atados"
$ws.Range("O35").Value = "Región Metropolitana"
$ws.Range("P35").Value = 667
$ws.Range("Q35").Value = 3
$ws.Range("R35").Value = "Hortaliza"
